# Apply updated symbol list values (generated Thu Jan  5 05:10:47 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices, percentages, hour) stay stored as text,
# matching the workbook's existing inline-string convention instead of being
# auto-coerced to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "258.60"
$ws.Range("E2").Value = "1.35%"
$ws.Range("G2").Value = "5"

# Row 3
$ws.Range("D3").Value = "26.88"
$ws.Range("E3").Value = "-4.03%"
$ws.Range("G3").Value = "5"

# Row 4
$ws.Range("D4").Value = "4.859"
$ws.Range("E4").Value = "-9.57%"
$ws.Range("G4").Value = "5"

# Row 5
$ws.Range("D5").Value = "0.05961"
$ws.Range("E5").Value = "2.27%"
$ws.Range("G5").Value = "5"

# Row 6
$ws.Range("D6").Value = "6.688"
$ws.Range("E6").Value = "-0.52%"
$ws.Range("G6").Value = "5"

# Row 7
$ws.Range("D7").Value = "0.8769"
$ws.Range("E7").Value = "0.79%"
$ws.Range("G7").Value = "5"

# Row 8
$ws.Range("D8").Value = "0.9582"
$ws.Range("E8").Value = "5.57%"
$ws.Range("G8").Value = "5"

# Row 9
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "0.0006077"
$ws.Range("E9").Value = "0.62%"
$ws.Range("G9").Value = "5"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1415"
$ws.Range("E10").Value = "-0.15%"
$ws.Range("G10").Value = "5"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07223"
$ws.Range("E11").Value = "-0.08%"
$ws.Range("G11").Value = "5"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.03147"
$ws.Range("E12").Value = "-1.00%"
$ws.Range("G12").Value = "5"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.09236"
$ws.Range("E13").Value = "-0.07%"
$ws.Range("G13").Value = "5"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001538"
$ws.Range("E14").Value = "-0.43%"
$ws.Range("G14").Value = "5"

# Row 15
$ws.Range("D15").Value = "0.005950"
$ws.Range("E15").Value = "-1.56%"
$ws.Range("G15").Value = "5"

# Row 16
$ws.Range("D16").Value = "3.484"
$ws.Range("E16").Value = "-0.25%"
$ws.Range("G16").Value = "5"

# Row 17
$ws.Range("D17").Value = "3.221"
$ws.Range("E17").Value = "-0.12%"
$ws.Range("G17").Value = "5"

# Row 18
$ws.Range("G18").Value = "5"

# Row 19
$ws.Range("D19").Value = "0.3144"
$ws.Range("E19").Value = "-0.77%"
$ws.Range("G19").Value = "5"

# Row 20
$ws.Range("D20").Value = "0.03602"
$ws.Range("E20").Value = "4.12%"
$ws.Range("G20").Value = "5"

# Row 21
$ws.Range("G21").Value = "5"

# Row 22
$ws.Range("D22").Value = "3.521"
$ws.Range("E22").Value = "-0.37%"
$ws.Range("G22").Value = "5"

# Row 23
$ws.Range("D23").Value = "0.04210"
$ws.Range("E23").Value = "1.50%"
$ws.Range("G23").Value = "5"

# Row 24
$ws.Range("E24").Value = "0.05%"
$ws.Range("G24").Value = "5"

# Row 25
$ws.Range("D25").Value = "0.001223"
$ws.Range("E25").Value = "-0.01%"
$ws.Range("G25").Value = "5"

# Row 26
$ws.Range("D26").Value = "0.004518"
$ws.Range("E26").Value = "-7.19%"
$ws.Range("G26").Value = "5"

# Row 27
$ws.Range("E27").Value = "-0.06%"
$ws.Range("G27").Value = "5"

# Row 28
$ws.Range("D28").Value = "0.0001492"
$ws.Range("E28").Value = "2.55%"
$ws.Range("G28").Value = "5"

# Row 29
$ws.Range("G29").Value = "5"

# Row 30
$ws.Range("G30").Value = "5"

# Row 31
$ws.Range("G31").Value = "5"

# Row 32
$ws.Range("G32").Value = "5"

# Row 33
$ws.Range("G33").Value = "5"

# Row 34
$ws.Range("G34").Value = "5"

# Row 35
$ws.Range("G35").Value = "5"

# Row 36
$ws.Range("G36").Value = "5"

# Row 37
$ws.Range("G37").Value = "5"

# Row 38
$ws.Range("G38").Value = "5"

# Row 39
$ws.Range("G39").Value = "5"

# Row 40
$ws.Range("D40").Value = "0.03838"
$ws.Range("E40").Value = "-0.30%"
$ws.Range("G40").Value = "5"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.005887"
$ws.Range("E41").Value = "2.31%"
$ws.Range("G41").Value = "5"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1104"
$ws.Range("E42").Value = "0.18%"
$ws.Range("G42").Value = "5"

# Row 43
$ws.Range("D43").Value = "0.001899"
$ws.Range("E43").Value = "-20.21%"
$ws.Range("G43").Value = "5"

# Row 44
$ws.Range("E44").Value = "5.43%"
$ws.Range("G44").Value = "5"

# Row 45
$ws.Range("D45").Value = "0.00005487"
$ws.Range("E45").Value = "3.99%"
$ws.Range("G45").Value = "5"

# Row 46
$ws.Range("E46").Value = "-0.04%"
$ws.Range("G46").Value = "5"

# Row 47
$ws.Range("D47").Value = "0.1090"
$ws.Range("E47").Value = "9.04%"
$ws.Range("G47").Value = "5"

# Row 48
$ws.Range("D48").Value = "0.002125"
$ws.Range("E48").Value = "-3.35%"
$ws.Range("G48").Value = "5"

# Row 49
$ws.Range("E49").Value = "-0.04%"
$ws.Range("G49").Value = "5"

# Row 50
$ws.Range("E50").Value = "-0.04%"
$ws.Range("G50").Value = "5"

# Row 51
$ws.Range("G51").Value = "5"
